$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("I4").Value = -1.304125163509162
$ws.Range("J4").Value = 0.4726427527162234
$ws.Range("K4").Value = 0.590360291702797
$ws.Range("L4").Value = 3.020339017269187
